$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 11:27:11"
$ws1.Cells.Item(3, 1).Value = "Total filas: 884"

# columns (B..G): Hora_Scrap, Hora_Llegada, Línea, Minutos, Parada, Fecha
$ws1Rows = @(872, 873, 874, 875, 876, 877, 878, 879, 880, 881, 882, 883, 884, 885)
$ws1ColB = @("11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00", "11:27:00")
$ws1ColC = @("11:30", "11:34", "11:35", "11:41", "11:46", "11:54", "11:55", "11:58", "12:04", "12:18", "12:19", "12:30", "12:34", "13:02")
$ws1ColD = @("11_ETCHEVERRY", "23_HERNANDEZ", "10_OLMOS", "215A_EL PATO", "16_SANTA ANA", "15_ABASTO", "225_GOMEZ", "16_SANTA ANA", "23_HERNANDEZ", "15_ABASTO", "10_OLMOS", "215C_EL PATO", "23_HERNANDEZ", "215C_EL PATO")
$ws1ColE = @(3, 7, 8, 14, 19, 27, 28, 31, 37, 51, 52, 63, 67, 95)
$ws1ColF = @("LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912", "LP1912")
$ws1ColG = @("31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025", "31/12/2025")
for ($i = 0; $i -lt $ws1Rows.Count; $i++) {
    $r = $ws1Rows[$i]
    $ws1.Cells.Item($r, 2).Value = $ws1ColB[$i]
    $ws1.Cells.Item($r, 3).Value = $ws1ColC[$i]
    $ws1.Cells.Item($r, 4).Value = $ws1ColD[$i]
    $ws1.Cells.Item($r, 5).Value = $ws1ColE[$i]
    $ws1.Cells.Item($r, 6).Value = $ws1ColF[$i]
    $ws1.Cells.Item($r, 7).Value = $ws1ColG[$i]
}

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 11:27:11"
$ws2.Cells.Item(3, 1).Value = "Total filas: 67"

# columns (B..G): Fecha, Hora_Scrap, Hora_Llegada, Línea, Minutos, Parada
$ws2Rows = @(66, 67, 68)
$ws2ColB = @("31/12/2025", "31/12/2025", "31/12/2025")
$ws2ColC = @("11:27:00", "11:27:00", "11:27:00")
$ws2ColD = @("11:41", "12:30", "13:02")
$ws2ColE = @("215A_EL PATO", "215C_EL PATO", "215C_EL PATO")
$ws2ColF = @(14, 63, 95)
$ws2ColG = @("LP1912", "LP1912", "LP1912")
for ($i = 0; $i -lt $ws2Rows.Count; $i++) {
    $r = $ws2Rows[$i]
    $ws2.Cells.Item($r, 2).Value = $ws2ColB[$i]
    $ws2.Cells.Item($r, 3).Value = $ws2ColC[$i]
    $ws2.Cells.Item($r, 4).Value = $ws2ColD[$i]
    $ws2.Cells.Item($r, 5).Value = $ws2ColE[$i]
    $ws2.Cells.Item($r, 6).Value = $ws2ColF[$i]
    $ws2.Cells.Item($r, 7).Value = $ws2ColG[$i]
}

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 11:27:11"
$ws3.Cells.Item(3, 1).Value = "Total filas: 103"

# columns (B..G): Fecha, Hora_Scrap, Hora_Llegada, Línea, Minutos, Parada
$ws3Rows = @(104)
$ws3ColB = @("31/12/2025")
$ws3ColC = @("11:27:06")
$ws3ColD = @("11:44")
$ws3ColE = @("215C_LA PLATA")
$ws3ColF = @(17)
$ws3ColG = @("L6203")
for ($i = 0; $i -lt $ws3Rows.Count; $i++) {
    $r = $ws3Rows[$i]
    $ws3.Cells.Item($r, 2).Value = $ws3ColB[$i]
    $ws3.Cells.Item($r, 3).Value = $ws3ColC[$i]
    $ws3.Cells.Item($r, 4).Value = $ws3ColD[$i]
    $ws3.Cells.Item($r, 5).Value = $ws3ColE[$i]
    $ws3.Cells.Item($r, 6).Value = $ws3ColF[$i]
    $ws3.Cells.Item($r, 7).Value = $ws3ColG[$i]
}
